# adattamento per filippine, in attesa di valori sigmaT
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("geo")

# Workbook calculation settings: switch to manual calculation and R1C1 reference mode
$excel.Calculation = -4135      # xlCalculationManual
$excel.ReferenceStyle = -4150   # xlR1C1

# Column G: multiply existing unit weight values by g = 9.81 (turn literals into formulas too)
$ws.Range("G4").Formula  = "=(2.75*0.8+2.1*0.2)*9.81"
$ws.Range("G5").Formula  = "=(2.75*0.8+2.1*0.2)*9.81"
$ws.Range("G6").Formula  = "=(2.75*0.8+2.1*0.2)*9.81"
$ws.Range("G7").Formula  = "=2.3*9.81"
$ws.Range("G8").Formula  = "=(2.75*0.8+2.1*0.2)*9.81"
$ws.Range("G9").Formula  = "=2.3*9.81"
$ws.Range("G10").Formula = "=(2.75*0.8+2.1*0.2)*9.81"
$ws.Range("G11").Formula = "=2.3*9.81"
$ws.Range("G12").Formula = "=(2.75*0.8+2.1*0.2)*9.81"
$ws.Range("G13").Formula = "=2.3*9.81"
$ws.Range("G14").Formula = "=2.1*9.81"
$ws.Range("G15").Formula = "=2.1*9.81"

# Column K: values expressed in kN instead of N (divide by 1000)
$ws.Range("K4").Value  = 0.3
$ws.Range("K5").Value  = 0.38
$ws.Range("K6").Value  = 0.4
$ws.Range("K7").Value  = 0.42
$ws.Range("K8").Value  = 0.4
$ws.Range("K9").Value  = 0.42
$ws.Range("K10").Value = 0.38
$ws.Range("K11").Value = 0.38
$ws.Range("K12").Value = 0.32
$ws.Range("K13").Value = 0.3
$ws.Range("K14").Value = 0.3
$ws.Range("K15").Value = 0.25

# Move the active selection down one row, from K15 to K16
$ws.Range("K16").Select()
